$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E are treated as Text so values like "65.532.61" are not
# reinterpreted as numbers/dates by Excel.
$ws.Columns.Item(4).NumberFormat = "@"
$ws.Columns.Item(5).NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "65.532.61"
$ws.Cells.Item(2, 5).Value = "  -2.79%  "
$ws.Cells.Item(3, 4).Value = "2.476.76"
$ws.Cells.Item(3, 5).Value = "  -5.76%  "
$ws.Cells.Item(4, 4).Value = "1.00"
$ws.Cells.Item(4, 5).Value = "  -0.01%  "
$ws.Cells.Item(5, 4).Value = "574.93"
$ws.Cells.Item(5, 5).Value = "  -3.40%  "
$ws.Cells.Item(6, 4).Value = "166.14"
$ws.Cells.Item(6, 5).Value = "  -0.97%  "
$ws.Cells.Item(7, 5).Value = "  +0.09%  "
$ws.Cells.Item(8, 4).Value = "0.515"
$ws.Cells.Item(8, 5).Value = "  -3.62%  "
$ws.Cells.Item(9, 4).Value = "2.473.65"
$ws.Cells.Item(9, 5).Value = "  -5.90%  "
$ws.Cells.Item(10, 4).Value = "0.134"
$ws.Cells.Item(10, 5).Value = "  -4.02%  "
$ws.Cells.Item(11, 5).Value = "  -0.51%  "
$ws.Cells.Item(12, 4).Value = "0.345"
$ws.Cells.Item(12, 5).Value = "  -5.59%  "
$ws.Cells.Item(13, 4).Value = "5.04"
$ws.Cells.Item(13, 5).Value = "  -3.80%  "
$ws.Cells.Item(14, 4).Value = "26.06"
$ws.Cells.Item(14, 5).Value = "  -5.97%  "
$ws.Cells.Item(15, 4).Value = "2.954.76"
$ws.Cells.Item(15, 5).Value = "  -4.88%  "
$ws.Cells.Item(16, 4).Value = "0.0000172"
$ws.Cells.Item(16, 5).Value = "  -6.08%  "
$ws.Cells.Item(17, 4).Value = "65.609.99"
$ws.Cells.Item(17, 5).Value = "  -2.41%  "
$ws.Cells.Item(18, 4).Value = "2.481.28"
$ws.Cells.Item(18, 5).Value = "  -5.02%  "
$ws.Cells.Item(19, 4).Value = "11.10"
$ws.Cells.Item(19, 5).Value = "  -8.31%  "
$ws.Cells.Item(20, 4).Value = "7.53"
$ws.Cells.Item(20, 5).Value = "  -6.03%  "
$ws.Cells.Item(21, 4).Value = "340.84"
$ws.Cells.Item(21, 5).Value = "  -5.02%  "
$ws.Cells.Item(22, 4).Value = "4.13"
$ws.Cells.Item(22, 5).Value = "  -4.59%  "
$ws.Cells.Item(23, 4).Value = "4.49"
$ws.Cells.Item(23, 5).Value = "  -4.28%  "
$ws.Cells.Item(24, 5).Value = "  +0.11%  "
$ws.Cells.Item(25, 4).Value = "1.91"
$ws.Cells.Item(25, 5).Value = "  -1.39%  "
$ws.Cells.Item(26, 4).Value = "68.44"
$ws.Cells.Item(26, 5).Value = "  -2.03%  "
$ws.Cells.Item(27, 4).Value = "9.74"
$ws.Cells.Item(27, 5).Value = "  -5.25%  "
$ws.Cells.Item(28, 4).Value = "0.999"
$ws.Cells.Item(28, 5).Value = "  -0.31%  "
$ws.Cells.Item(29, 4).Value = "2.605.99"
$ws.Cells.Item(29, 5).Value = "  -5.57%  "
$ws.Cells.Item(30, 4).Value = "0.0₃0951"
$ws.Cells.Item(30, 5).Value = "  -5.52%  "
$ws.Cells.Item(31, 4).Value = "514.14"
$ws.Cells.Item(31, 5).Value = "  -5.49%  "
$ws.Cells.Item(32, 4).Value = "7.99"
$ws.Cells.Item(32, 5).Value = "  +0.51%  "
$ws.Cells.Item(33, 4).Value = "1.29"
$ws.Cells.Item(33, 5).Value = "  -4.49%  "
$ws.Cells.Item(34, 4).Value = "1.79"
$ws.Cells.Item(34, 5).Value = "  -5.81%  "
$ws.Cells.Item(35, 4).Value = "0.129"
$ws.Cells.Item(35, 5).Value = "  -5.00%  "
$ws.Cells.Item(36, 4).Value = "0.998"
$ws.Cells.Item(36, 5).Value = "  -0.22%  "
$ws.Cells.Item(37, 4).Value = "156.71"
$ws.Cells.Item(37, 5).Value = "  +0.03%  "
$ws.Cells.Item(38, 4).Value = "1.43"
$ws.Cells.Item(38, 5).Value = "  -5.19%  "
$ws.Cells.Item(39, 4).Value = "18.30"
$ws.Cells.Item(39, 5).Value = "  -3.92%  "
$ws.Cells.Item(40, 4).Value = "18.25"
$ws.Cells.Item(40, 5).Value = "  +0.22%  "
$ws.Cells.Item(41, 4).Value = "0.349"
$ws.Cells.Item(41, 5).Value = "  -4.82%  "
$ws.Cells.Item(42, 4).Value = "1.74"
$ws.Cells.Item(42, 5).Value = "  -4.11%  "
$ws.Cells.Item(43, 4).Value = "4.97"
$ws.Cells.Item(43, 5).Value = "  -4.80%  "
$ws.Cells.Item(44, 5).Value = "  -0.02%  "
$ws.Cells.Item(45, 4).Value = "2.37"
$ws.Cells.Item(45, 5).Value = "  -2.57%  "
$ws.Cells.Item(46, 4).Value = "144.84"
$ws.Cells.Item(46, 5).Value = "  -4.99%  "
$ws.Cells.Item(47, 4).Value = "0.547"
$ws.Cells.Item(47, 5).Value = "  -5.83%  "
$ws.Cells.Item(48, 4).Value = "3.63"
$ws.Cells.Item(48, 5).Value = "  -4.24%  "
$ws.Cells.Item(49, 4).Value = "0.0₆0266"
$ws.Cells.Item(49, 5).Value = "  -10.60%  "
$ws.Cells.Item(50, 4).Value = "1.67"
$ws.Cells.Item(50, 5).Value = "  -1.95%  "
$ws.Cells.Item(51, 4).Value = "0.0744"
$ws.Cells.Item(51, 5).Value = "  -3.47%  "
